# Update data values in the KNN imputation result sheet.
# The underlying data table (header row 1: A,B,C,D,E ; data rows 2-102)
# has several of its numeric values updated to reflect new algorithm output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.64
$ws.Range("C3").Value = -12.774
$ws.Range("D19").Value = -8.001999999999999
$ws.Range("A21").Value = -20.257
$ws.Range("A23").Value = -20.955
$ws.Range("C24").Value = -12.658
$ws.Range("D24").Value = -7.657999999999999
$ws.Range("A25").Value = -21.671
$ws.Range("B27").Value = 5.513
$ws.Range("D30").Value = -7.234999999999999
$ws.Range("B31").Value = 6.016
$ws.Range("D31").Value = -8.127000000000001
$ws.Range("D33").Value = -7.622
$ws.Range("B39").Value = 7.915000000000001
$ws.Range("B48").Value = 5.274
$ws.Range("B51").Value = 5.962
$ws.Range("B52").Value = 5.295
$ws.Range("A53").Value = -22.01
$ws.Range("B55").Value = 4.508
$ws.Range("D55").Value = -8.221
$ws.Range("B56").Value = 5.169
$ws.Range("A57").Value = -21.768
$ws.Range("B57").Value = 6.531000000000001
$ws.Range("C57").Value = -12.98
$ws.Range("A59").Value = -22.363
$ws.Range("C61").Value = -13.508
$ws.Range("D65").Value = -7.869
$ws.Range("A69").Value = -21.634
$ws.Range("C70").Value = -12.166
$ws.Range("D70").Value = -7.5
$ws.Range("B73").Value = 7.343000000000001
$ws.Range("D75").Value = -7.702
$ws.Range("A79").Value = -21.115
$ws.Range("A83").Value = -22.134
$ws.Range("D83").Value = -8.512
$ws.Range("C86").Value = -13.597
$ws.Range("B89").Value = 5.999
$ws.Range("B90").Value = 5.811
$ws.Range("A93").Value = -21.536
$ws.Range("D96").Value = -7.411000000000001
$ws.Range("D97").Value = -8.166
$ws.Range("C98").Value = -12.45
$ws.Range("C100").Value = -12.621
$ws.Range("C102").Value = -13.564
